$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text needs to change. Every one of these values is stored as
# TEXT (shared string), even though it looks numeric, so we briefly force a
# text number-format before writing the value (otherwise Excel infers a
# number) and then clear the format again so no residual cell formatting is
# left behind.
$targets = @("B6", "B11", "B12", "B13", "D11", "D12", "D13")
foreach ($addr in $targets) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("B6").Value  = "0.26607752571579923"
$ws.Range("B11").Value = "0.617737003058104"
$ws.Range("B12").Value = "1"
$ws.Range("B13").Value = "0.0"
$ws.Range("D11").Value = "0.6183206106870229"
$ws.Range("D12").Value = "1"
$ws.Range("D13").Value = "0.0"

foreach ($addr in $targets) {
    $ws.Range($addr).ClearFormats()
}
